$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '58.120.65'
$ws.Range('E2').Value = '  -1.33%  '
$ws.Range('D3').Value = '2.476.78'
$ws.Range('E3').Value = '  -1.54%  '
$ws.Range('E4').Value = '  +0.35%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '519.60'
$ws.Range('E5').Value = '  -2.96%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '131.20'
$ws.Range('E6').Value = '  -3.99%  '
$ws.Range('E8').Value = '  -1.59%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.0993'
$ws.Range('E9').Value = '  -1.36%  '
$ws.Range('E10').Value = '  -0.43%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.33'
$ws.Range('E11').Value = '  +0.25%  '
$ws.Range('E12').Value = '  -0.57%  '
$ws.Range('D13').Value = '2.917.29'
$ws.Range('E13').Value = '  -1.38%  '
$ws.Range('D14').Value = '58.068.30'
$ws.Range('E14').Value = '  -1.30%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '22.35'
$ws.Range('E15').Value = '  -2.33%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000136'
$ws.Range('E16').Value = '  -1.60%  '
$ws.Range('D17').Value = '2.482.29'
$ws.Range('E17').Value = '  -0.96%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '10.84'
$ws.Range('E18').Value = '  -2.01%  '
$ws.Range('B19').Value = 'BitcoinCash'
$ws.Range('C19').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '320.64'
$ws.Range('E19').Value = '  -0.67%  '
$ws.Range('B20').Value = 'Polkadot'
$ws.Range('C20').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.18'
$ws.Range('E20').Value = '  -2.00%  '
$ws.Range('E21').Value = '  +0.00%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.78'
$ws.Range('E22').Value = '  -2.64%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '64.22'
$ws.Range('E23').Value = '  -2.22%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.410'
$ws.Range('E24').Value = '  -2.21%  '
$ws.Range('E25').Value = '  +0.25%  '
$ws.Range('E26').Value = '  -2.51%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.34'
$ws.Range('E27').Value = '  -2.26%  '
$ws.Range('D28').Value = '0.0₃0755'
$ws.Range('E28').Value = '  -1.79%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.70'
$ws.Range('E29').Value = '  -3.71%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.31'
$ws.Range('E30').Value = '  -5.58%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '165.83'
$ws.Range('E31').Value = '  -0.15%  '
$ws.Range('E32').Value = '  +0.48%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.998'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.999'
$ws.Range('E34').Value = '  -0.03%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '18.10'
$ws.Range('E35').Value = '  -1.70%  '
$ws.Range('E36').Value = '  -9.04%  '
$ws.Range('E37').Value = '  -2.18%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.47'
$ws.Range('E38').Value = '  -3.80%  '
$ws.Range('E39').Value = '  -2.40%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.47'
$ws.Range('E40').Value = '  -3.50%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '275.58'
$ws.Range('E41').Value = '  -2.82%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.02'
$ws.Range('E42').Value = '  -3.78%  '
$ws.Range('E43').Value = '  -1.61%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '127.21'
$ws.Range('E44').Value = '  -3.94%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0907'
$ws.Range('E45').Value = '  -1.79%  '
$ws.Range('E46').Value = '  -3.11%  '
$ws.Range('E47').Value = '  -2.40%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '17.14'
$ws.Range('E48').Value = '  -0.42%  '
$ws.Range('D49').Value = '1.742.65'
$ws.Range('E49').Value = '  -1.34%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.973'
$ws.Range('E50').Value = '  -0.99%  '
$ws.Range('E51').Value = '  -1.05%  '
